# Fixed update to excel issue
#
# 1. Rename "Requested quantity" header on "Weekly Quantity" -> "Weekly_PO_Qty"
# 2. Rename "Requested quantity" header on "Monthly Trend"   -> "Monthly_PO_Qty"
# 3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

# --- 1 & 2: header renames -------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: new "PO Forecast" sheet --------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Headers - reuse the same header style (bold/centered/bordered) that is
# already used on "Weekly Quantity"!A1:B1 by copying formats across, then
# set the header text.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45095.99999999999, 3,  1.019439895130605, 5.273635334901883),
    @(45116.99999999999, 5,  3.386839873195119, 7.48906471991998),
    @(45130.99999999999, 7,  4.908201414285021, 9.062317138994654),
    @(45144.99999999999, 9,  6.35519594017163,  10.62754544422196),
    @(45151.99999999999, 9,  7.16474199560646,  11.37436149934012),
    @(45165.99999999999, 11, 8.813594264377459, 12.80086729355572),
    @(45172.99999999999, 12, 9.445765813217941, 13.61184321768039),
    @(45179.99999999999, 12, 10.25609941528958, 14.4167662395687),
    @(45186.99999999999, 13, 11.05875956799991, 15.29233047993442),
    @(45193.99999999999, 14, 11.82749428144212, 15.8787952070415),
    @(45200.99999999999, 15, 12.76894936374831, 16.78105933834741),
    @(45207.99999999999, 15, 13.37557166289719, 17.44563968084685),
    @(45214.99999999999, 16, 14.29847210622072, 18.31657343760082),
    @(45221.99999999999, 17, 14.92787014739415, 19.06821368718407),
    @(45228.99999999999, 18, 15.7550355560672,  19.79223921119289)
)

$row = 2
foreach ($entry in $data) {
    $dsCell = $wsForecast.Cells.Item($row, 1)
    $dsCell.Value = $entry[0]
    $dsCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $wsForecast.Cells.Item($row, 2).Value = $entry[1]
    $wsForecast.Cells.Item($row, 3).Value = $entry[2]
    $wsForecast.Cells.Item($row, 4).Value = $entry[3]

    $row++
}
